# Applies the Tanmenet.xlsx edit described by the commit:
# "Add initial documentation for container concepts and usage"
#
# Concretely this:
#  - adds a "Flappy bird player" entry (C6) next to the existing
#    "Flappy bird clouds" entry (C5) on the "Emelt" sheet
#  - marks two previously-scheduled "Méh" lessons (D7, D9) as cancelled,
#    replacing them with "ELMARADT" notes in column C (C7, C9)
#  - adds an "Adatbázis.ppt" note in C8
#  - updates the selection / cursor position left behind in the sheet
#  - normalizes the header/footer font name from "Általános" to "Regular"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Emelt")
$ws2 = $wb.Worksheets.Item("Közép")

# --- New "Flappy bird player" row, alongside existing "Flappy bird clouds" ---
$ws.Range("C6").Value = "Flappy bird player"

# --- Row 7: lesson "ELMARADT" (cancelled); clear the old "Méh" marker ---
$ws.Range("C7").Value = "ELMARADT"
$ws.Range("D7").ClearContents()

# --- Row 8: add "Adatbázis.ppt" note, keep existing "Méh" marker in D8 ---
$ws.Range("C8").Value = "Adatbázis.ppt"

# --- Row 9: lesson "ELMARADT" (cancelled); clear the old "Méh" marker ---
$ws.Range("C9").Value = "ELMARADT"
$ws.Range("D9").ClearContents()

# --- Leave the cursor on D9, matching the saved selection ---
$ws.Range("D9").Select() | Out-Null

# --- Normalize header/footer font name on both sheets ---
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'

$ws2.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws2.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'
